# Testing commit on 07062021
#
# Updates the sample login-data sheet:
#   - Row 1 credentials/url get swapped out for a new test login
#     (admin -> nilanjan / admin@123 -> Admin@123 / url gets a deep link)
#   - Leaves rows 2-4 (sathish/manikanta/rao1232 + admin@123 + Pass) untouched
#   - Re-positions the workbook window and the active selection, matching the
#     state the workbook was left in after the edit was made in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new login under test --------------------------------------
# Write G1/B1 before A1 so the shared-string table is appended in the same
# order the source workbook shows (url, password, username).
$ws.Range("G1").Value = "https://ostaffuidev.onpassive.com/login/loginPage"
$ws.Range("B1").Value = "Admin@123"
$ws.Range("A1").Value = "nilanjan"

# --- Window geometry (bookViews/workbookView) ---------------------------
try {
    $win = $excel.Windows.Item(1)
    $win.Width  = 15540
    $win.Height = 3165
    $win.Left   = 4320
    $win.Top    = 3105
} catch {
    # window-position persistence is a cosmetic, best-effort change only
}

# --- Selection: A2 selected first, then A1 becomes the active cell -----
try {
    $ws.Range("A2").Select() | Out-Null
    $ws.Range("A1").Select() | Out-Null
} catch {
    # selection state is cosmetic; ignore if unsupported
}
